$wb = $excel.ActiveWorkbook

function Add-ResultSheet {
    param(
        [string]$SheetName,
        [string]$Label,
        [double]$InVehicle,
        [double]$AtStop,
        [double]$Extra,
        [double]$Tardiness,
        [double]$Total
    )

    # Insert the new sheet right after the current last sheet, so it lands
    # at the end of the tab list (matching the order sheets are appended
    # in the workbook).
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $SheetName

    # Header row (B1:F1) - bold, thin border box, centered horizontally,
    # top-aligned vertically (matches the style used on every other sheet
    # in this workbook).
    $headerRange = $ws.Range("B1:F1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
    $headerRange.Borders.LineStyle = 1
    $headerRange.Borders.Weight = 2

    $ws.Range("B1").Value = "In-vehicle"
    $ws.Range("C1").Value = "At-stop"
    $ws.Range("D1").Value = "Extra"
    $ws.Range("E1").Value = "Tardiness"
    $ws.Range("F1").Value = "Total"

    # Row label (A2) gets the same styling as the header cells.
    $a2 = $ws.Range("A2")
    $a2.Font.Bold = $true
    $a2.HorizontalAlignment = -4108
    $a2.VerticalAlignment = -4160
    $a2.Borders.LineStyle = 1
    $a2.Borders.Weight = 2
    $a2.Value = $Label

    $ws.Range("B2").Value = $InVehicle
    $ws.Range("C2").Value = $AtStop
    $ws.Range("D2").Value = $Extra
    $ws.Range("E2").Value = $Tardiness
    $ws.Range("F2").Value = $Total
}

Add-ResultSheet "FTNC_Demand516" "FTNC" 2109.323368643516 12732.99064835311 374.5119430131296 43.50906907283795 15260.33502908257
Add-ResultSheet "FTHC_Demand51" "FTHC" 2232.164806829644 12581.78707840258 375.4875219215975 5.412815180977432 15194.85222233474

# Restore the originally active sheet/tab selection (sheet creation shifts
# focus onto the newly added sheets).
$wb.Worksheets.Item(1).Activate()
